$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new columns I (I0) and J (IF), matching style of H1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null

# Data for columns I and J, rows 2-39
$data = @{
    2  = @(1, 5)
    3  = @(1, 4)
    4  = @(1, 8)
    5  = @(1, 6)
    6  = @(1, 5)
    7  = @(1, 5)
    8  = @(1, 6)
    9  = @(1, 6)
    10 = @(1, 6)
    11 = @(1, 5)
    12 = @(1, 6)
    13 = @(1, 6)
    14 = @(1, 7)
    15 = @(1, 4)
    16 = @(1, 6)
    17 = @(1, 8)
    18 = @(1, 6)
    19 = @(1, 6)
    20 = @(1, 5)
    21 = @(1, 6)
    22 = @(1, 6)
    23 = @(1, 6)
    24 = @(1, 6)
    25 = @(1, 7)
    26 = @(1, 5)
    27 = @(1, 7)
    28 = @(1, 5)
    29 = @(1, 6)
    30 = @(1, 5)
    31 = @(1, 6)
    32 = @(1, 6)
    33 = @(1, 6)
    34 = @(1, 5)
    35 = @(1, 6)
    36 = @(1, 5)
    37 = @(1, 4)
    38 = @(6, 7)
    39 = @(3, 4)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 9).Value = $vals[0]
    $ws.Cells.Item($r, 10).Value = $vals[1]
}
